$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AssessmentResult")

# Update E16 value
$ws.Range("E16").Value = 2

# Update row 40 values (B40:H40)
$ws.Range("B40").Value = 3
$ws.Range("C40").Value = 2
$ws.Range("D40").Value = 2
$ws.Range("E40").Value = 4
$ws.Range("F40").Value = 2
$ws.Range("G40").Value = 3
$ws.Range("H40").Value = 3

# Update the active selection on the sheet to G13
$ws.Activate()
$ws.Range("G13").Select()
